# Update the two-digit-divided-by-one-digit division worksheet numbers.
# The document contains a single table; the division problems live in
# rows 1, 5, 9, 13, 17 (1-based), five per row. Each cell's value is
# replaced in place via the Tables object model so that run formatting
# (font/size) and paragraph structure are preserved untouched.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "34÷5="
$t.Cell(1,2).Range.Text  = "85÷2="
$t.Cell(1,3).Range.Text  = "60÷6="
$t.Cell(1,4).Range.Text  = "42÷7="
$t.Cell(1,5).Range.Text  = "17÷9="

$t.Cell(5,1).Range.Text  = "39÷2="
$t.Cell(5,2).Range.Text  = "99÷7="
$t.Cell(5,3).Range.Text  = "71÷7="
$t.Cell(5,4).Range.Text  = "46÷7="
$t.Cell(5,5).Range.Text  = "51÷2="

$t.Cell(9,1).Range.Text  = "59÷2="
$t.Cell(9,2).Range.Text  = "96÷7="
$t.Cell(9,3).Range.Text  = "75÷4="
$t.Cell(9,4).Range.Text  = "29÷3="
$t.Cell(9,5).Range.Text  = "80÷8="

$t.Cell(13,1).Range.Text = "57÷9="
$t.Cell(13,2).Range.Text = "20÷8="
$t.Cell(13,3).Range.Text = "78÷7="
$t.Cell(13,4).Range.Text = "73÷3="
$t.Cell(13,5).Range.Text = "97÷5="

$t.Cell(17,1).Range.Text = "27÷3="
$t.Cell(17,2).Range.Text = "29÷6="
$t.Cell(17,3).Range.Text = "80÷6="
$t.Cell(17,4).Range.Text = "86÷3="
$t.Cell(17,5).Range.Text = "30÷4="
